$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.440.73'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '3.514.63'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.78%  '
$ws.Range("E7").Value = '  -1.03%  '
$ws.Range("D8").Value = '3.506.22'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E10").Value = '  -3.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.583'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.25'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000275'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").Value = '4.074.90'
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '611.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '3.523.72'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").Value = '70.553.18'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.877'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '98.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.40%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  -1.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("E31").Value = '  -3.39%  '
$ws.Range("E32").Value = '  -4.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '635.88'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.87%  '
$ws.Range("E34").Value = '  -4.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.82'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.58'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.34%  '
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0475'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.142'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.72%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0743'
$ws.Range("E43").Value = '  +6.07%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.371.42'
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.308'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.36%  '
$ws.Range("E46").Value = '  -2.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("E49").Value = '  +0.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.64%  '
